$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data used to have a single header row (code/col1/col2/value3).
# It now represents a pandas MultiIndex on the columns, so a new top
# header row (level-0 labels "колонки"/"колонка") is added above the
# existing header row, and a blank separator row is left between the
# header block and the data rows.

# 1) Insert a brand-new row 1 for the level-0 MultiIndex labels.
#    This pushes the old header row (code/col1/col2/value3) down to row 2
#    and the data rows (BLR/KAZ/RUS/UKR) down to rows 3-6.
$ws.Rows.Item(1).Insert()

# 2) Insert a blank separator row right after the header (new row 3),
#    pushing the data rows down to their final positions (rows 4-7).
$ws.Rows.Item(3).Insert()
$ws.Rows.Item(3).ClearFormats()

# 3) Merge B1:C1 while the row is still unformatted, so Excel doesn't
#    need to split the shared border style between the two halves.
$ws.Range("B1:C1").Merge()

# 4) Copy the header row's formatting (bold font, border, centered/top
#    aligned) down onto the new row 1.
$ws.Range("A2:D2").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)

# 5) Fill in the level-0 MultiIndex labels.
$ws.Range("B1").Value = "колонки"
$ws.Range("D1").Value = "колонка"
